# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice[/NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ])
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 66.333336
$ws.Range("I4").Value = 55
$ws.Range("K4").Value = 55
$ws.Range("M4").Value = 59

$ws.Range("H62").Value = 4995
$ws.Range("I62").Value = 4995
$ws.Range("K62").Value = 4995
$ws.Range("M62").Value = -4371

$ws.Range("H65").Value = 4995
$ws.Range("I65").Value = 4995
$ws.Range("K65").Value = 24975
$ws.Range("M65").Value = -21855

$ws.Range("H127").Value = 8075.4
$ws.Range("I127").Value = 2804.7144
$ws.Range("K127").Value = 8414.143199999999
$ws.Range("M127").Value = -3454.143199999999

$ws.Range("H129").Value = 184996.55
$ws.Range("J129").Value = 2496
$ws.Range("L129").Value = 7488
$ws.Range("N129").Value = -17488

$ws.Range("H132").Value = 20708
$ws.Range("I132").Value = 3128.7727
$ws.Range("K132").Value = 9386.3181
$ws.Range("M132").Value = -6856.3181

$ws.Range("H139").Value = 51111.43
$ws.Range("J139").Value = 51111.43
$ws.Range("L139").Value = 51111.43
$ws.Range("N139").Value = -61391.43

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2233.9678
$ws.Range("I132").Value = 2227.4644
$ws.Range("K132").Value = 6682.3932
$ws.Range("M132").Value = -4152.3932

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1101
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 1101
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 1101
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -1551

$ws.Range("H67").Value = 1101
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 1101
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 1101
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -2661

$ws.Range("H105").Value = 2099.4
$ws.Range("I105").Value = 874.4375
$ws.Range("J105").Value = 4277.1113
$ws.Range("K105").Value = 874.4375
$ws.Range("L105").Value = 4277.1113
$ws.Range("M105").Value = 872.5625
$ws.Range("N105").Value = -7771.1113

$ws.Range("H134").Value = 2484.348
$ws.Range("J134").Value = 2806.1428
$ws.Range("L134").Value = 8418.428400000001
$ws.Range("N134").Value = -13488.4284

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1738.5238
$ws.Range("I31").Value = 1625.45
$ws.Range("J31").Value = 4000
$ws.Range("K31").Value = 1625.45
$ws.Range("L31").Value = 4000
$ws.Range("M31").Value = -1330.45
$ws.Range("N31").Value = -4590

$ws.Range("H34").Value = 1738.5238
$ws.Range("I34").Value = 1625.45
$ws.Range("J34").Value = 4000
$ws.Range("K34").Value = 1625.45
$ws.Range("L34").Value = 4000
$ws.Range("M34").Value = -1423.45
$ws.Range("N34").Value = -4404

$ws.Range("H99").Value = 9758975
$ws.Range("I99").Value = 1628961.8
$ws.Range("K99").Value = 1628961.8
$ws.Range("M99").Value = -1627463.8

$ws.Range("H122").Value = 235681.9
$ws.Range("I122").Value = 321795.38
$ws.Range("K122").Value = 965386.14
$ws.Range("M122").Value = -962936.14

$ws.Range("H126").Value = 9758975
$ws.Range("I126").Value = 1628961.8
$ws.Range("K126").Value = 4886885.4
$ws.Range("M126").Value = -4884415.4

$ws.Range("H132").Value = 5675.7144
$ws.Range("I132").Value = 3183.25
$ws.Range("K132").Value = 9549.75
$ws.Range("M132").Value = -7019.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 73.5
$ws.Range("I11").Value = 55.5
$ws.Range("J11").Value = 109.5
$ws.Range("K11").Value = 166.5
$ws.Range("L11").Value = 328.5
$ws.Range("M11").Value = -26.5
$ws.Range("N11").Value = -608.5

$ws.Range("H23").Value = 2176.6667
$ws.Range("J23").Value = 999.6667
$ws.Range("L23").Value = 2999.0001
$ws.Range("N23").Value = -3469.0001

$ws.Range("H44").Value = 55555624
$ws.Range("J44").Value = 83333390
$ws.Range("L44").Value = 250000170
$ws.Range("N44").Value = -250000966

$ws.Range("H54").Value = 75574.836
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 75574.836
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 226724.508
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -227842.508

$ws.Range("H60").Value = 875.3333
$ws.Range("I60").Value = 195.54546
$ws.Range("J60").Value = 2744.75
$ws.Range("K60").Value = 586.6363799999999
$ws.Range("L60").Value = 8234.25
$ws.Range("M60").Value = -335.6363799999999
$ws.Range("N60").Value = -8736.25

$ws.Range("H113").Value = 1749.8182
$ws.Range("I113").Value = 1935.3334
$ws.Range("J113").Value = 1527.2
$ws.Range("K113").Value = 5806.0002
$ws.Range("L113").Value = 4581.6
$ws.Range("M113").Value = -3636.0002
$ws.Range("N113").Value = -8921.6

$ws.Range("H131").Value = 3315.4119
$ws.Range("I131").Value = 1594.25
$ws.Range("J131").Value = 4845.3335
$ws.Range("K131").Value = 4782.75
$ws.Range("L131").Value = 14536.0005
$ws.Range("M131").Value = 257.25
$ws.Range("N131").Value = -24616.0005

$ws.Range("H139").Value = 333339330
$ws.Range("I139").Value = 1000000000
$ws.Range("K139").Value = 3000000000
$ws.Range("M139").Value = -2999994860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 6783.5884
$ws.Range("J113").Value = 8680
$ws.Range("L113").Value = 8680
$ws.Range("N113").Value = -13020

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1785
$ws.Range("I40").Value = 1785
$ws.Range("K40").Value = 1785
$ws.Range("M40").Value = -1649

$ws.Range("H100").Value = 6205.4707
$ws.Range("I100").Value = 4197.6665
$ws.Range("J100").Value = 7300.636
$ws.Range("K100").Value = 4197.6665
$ws.Range("L100").Value = 7300.636
$ws.Range("M100").Value = -3656.6665
$ws.Range("N100").Value = -8382.636

$ws.Range("H119").Value = 41210
$ws.Range("J119").Value = 41210
$ws.Range("L119").Value = 41210
$ws.Range("N119").Value = -50886

$ws.Range("H122").Value = 2790.5789
$ws.Range("I122").Value = 2725.6
$ws.Range("J122").Value = 3034.25
$ws.Range("K122").Value = 8176.799999999999
$ws.Range("L122").Value = 9102.75
$ws.Range("M122").Value = -5726.799999999999
$ws.Range("N122").Value = -14002.75

$ws.Range("H132").Value = 2626.6
$ws.Range("I132").Value = 2338.2
$ws.Range("K132").Value = 7014.599999999999
$ws.Range("M132").Value = -4484.599999999999

$ws.Range("H136").Value = 3397.44
$ws.Range("I136").Value = 2981.8235
$ws.Range("K136").Value = 8945.470499999999
$ws.Range("M136").Value = -6395.470499999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2443.6667
$ws.Range("I122").Value = 1666.5
$ws.Range("J122").Value = 3998
$ws.Range("K122").Value = 4999.5
$ws.Range("L122").Value = 11994
$ws.Range("M122").Value = -2549.5
$ws.Range("N122").Value = -16894

$ws.Range("H132").Value = 4168.968
$ws.Range("I132").Value = 3824
$ws.Range("J132").Value = 5012.222
$ws.Range("K132").Value = 11472
$ws.Range("L132").Value = 15036.666
$ws.Range("M132").Value = -8942
$ws.Range("N132").Value = -20096.666
